# edit.ps1
# Adds a new "2022-Q3" sheet (holdings-by-fund detail) right after "总计",
# and inserts a corresponding summary row at the top of the "总计" sheet's data.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q3" detail sheet, positioned before "2022-Q2"
#    (i.e. right after "总计"), mirroring the layout of the other
#    quarterly fund-holding sheets.
# ------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

# Header row (B1:H1)
$hdrRange = $wsQ3.Range("B1:H1")
$hdrRange.NumberFormat = "@"
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"
$hdrRange.Font.Bold = $true
$hdrRange.Borders.LineStyle = 1
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4160

# Data rows 2-8 (one row per fund). Columns A (index) and H (rank) are
# numeric; B-G mirror the source formatting (stored as text).
$rowRange = $wsQ3.Range("B2:G2")
$rowRange.NumberFormat = "@"
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("A2").Font.Bold = $true
$wsQ3.Range("A2").Borders.LineStyle = 1
$wsQ3.Range("A2").HorizontalAlignment = -4108
$wsQ3.Range("A2").VerticalAlignment = -4160
$wsQ3.Range("B2").Value = "009693"
$wsQ3.Range("C2").Value = "富国积极成长一年定期开放混合"
$wsQ3.Range("D2").Value = "12.30"
$wsQ3.Range("E2").Value = "97.80"
$wsQ3.Range("F2").Value = "3.81"
$wsQ3.Range("G2").Value = "0.4686"
$wsQ3.Range("H2").Value = 7

$rowRange = $wsQ3.Range("B3:G3")
$rowRange.NumberFormat = "@"
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("A3").Font.Bold = $true
$wsQ3.Range("A3").Borders.LineStyle = 1
$wsQ3.Range("A3").HorizontalAlignment = -4108
$wsQ3.Range("A3").VerticalAlignment = -4160
$wsQ3.Range("B3").Value = "014273"
$wsQ3.Range("C3").Value = "广发北交所精选两年定开混合A"
$wsQ3.Range("D3").Value = "3.37"
$wsQ3.Range("E3").Value = "64.25"
$wsQ3.Range("F3").Value = "4.62"
$wsQ3.Range("G3").Value = "0.1557"
$wsQ3.Range("H3").Value = 7

$rowRange = $wsQ3.Range("B4:G4")
$rowRange.NumberFormat = "@"
$wsQ3.Range("A4").Value = 2
$wsQ3.Range("A4").Font.Bold = $true
$wsQ3.Range("A4").Borders.LineStyle = 1
$wsQ3.Range("A4").HorizontalAlignment = -4108
$wsQ3.Range("A4").VerticalAlignment = -4160
$wsQ3.Range("B4").Value = "014269"
$wsQ3.Range("C4").Value = "嘉实北交所精选两年定期混合A"
$wsQ3.Range("D4").Value = "2.72"
$wsQ3.Range("E4").Value = "90.37"
$wsQ3.Range("F4").Value = "4.80"
$wsQ3.Range("G4").Value = "0.1306"
$wsQ3.Range("H4").Value = 9

$rowRange = $wsQ3.Range("B5:G5")
$rowRange.NumberFormat = "@"
$wsQ3.Range("A5").Value = 3
$wsQ3.Range("A5").Font.Bold = $true
$wsQ3.Range("A5").Borders.LineStyle = 1
$wsQ3.Range("A5").HorizontalAlignment = -4108
$wsQ3.Range("A5").VerticalAlignment = -4160
$wsQ3.Range("B5").Value = "014275"
$wsQ3.Range("C5").Value = "易方达北交所精选两年定开混合A"
$wsQ3.Range("D5").Value = "3.58"
$wsQ3.Range("E5").Value = "61.75"
$wsQ3.Range("F5").Value = "2.49"
$wsQ3.Range("G5").Value = "0.0891"
$wsQ3.Range("H5").Value = 8

$rowRange = $wsQ3.Range("B6:G6")
$rowRange.NumberFormat = "@"
$wsQ3.Range("A6").Value = 4
$wsQ3.Range("A6").Font.Bold = $true
$wsQ3.Range("A6").Borders.LineStyle = 1
$wsQ3.Range("A6").HorizontalAlignment = -4108
$wsQ3.Range("A6").VerticalAlignment = -4160
$wsQ3.Range("B6").Value = "014274"
$wsQ3.Range("C6").Value = "广发北交所精选两年定开混合C"
$wsQ3.Range("D6").Value = "0.85"
$wsQ3.Range("E6").Value = "64.25"
$wsQ3.Range("F6").Value = "4.62"
$wsQ3.Range("G6").Value = "0.0393"
$wsQ3.Range("H6").Value = 7

$rowRange = $wsQ3.Range("B7:G7")
$rowRange.NumberFormat = "@"
$wsQ3.Range("A7").Value = 5
$wsQ3.Range("A7").Font.Bold = $true
$wsQ3.Range("A7").Borders.LineStyle = 1
$wsQ3.Range("A7").HorizontalAlignment = -4108
$wsQ3.Range("A7").VerticalAlignment = -4160
$wsQ3.Range("B7").Value = "014270"
$wsQ3.Range("C7").Value = "嘉实北交所精选两年定期混合C"
$wsQ3.Range("D7").Value = "0.53"
$wsQ3.Range("E7").Value = "90.37"
$wsQ3.Range("F7").Value = "4.80"
$wsQ3.Range("G7").Value = "0.0254"
$wsQ3.Range("H7").Value = 9

$rowRange = $wsQ3.Range("B8:G8")
$rowRange.NumberFormat = "@"
$wsQ3.Range("A8").Value = 6
$wsQ3.Range("A8").Font.Bold = $true
$wsQ3.Range("A8").Borders.LineStyle = 1
$wsQ3.Range("A8").HorizontalAlignment = -4108
$wsQ3.Range("A8").VerticalAlignment = -4160
$wsQ3.Range("B8").Value = "014276"
$wsQ3.Range("C8").Value = "易方达北交所精选两年定开混合C"
$wsQ3.Range("D8").Value = "0.92"
$wsQ3.Range("E8").Value = "61.75"
$wsQ3.Range("F8").Value = "2.49"
$wsQ3.Range("G8").Value = "0.0229"
$wsQ3.Range("H8").Value = 8

# ------------------------------------------------------------------
# 2) Insert the matching summary row into "总计", above the existing
#    "2022-Q2"/"2022-Q1" rows, and renumber the index column (A).
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows(2).Insert()

$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 0.93

# Renumber the shifted rows (formerly "2022-Q2" -> row3, "2022-Q1" -> row4)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# Match A2 styling (bold, bordered, centered) to the other index cells
$a2 = $wsTotal.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

